$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt 12; $i++) {
    $srcRow = 2 + $i
    $destRow = 26 + $i
    for ($col = 1; $col -le 12; $col++) {
        $val = $ws.Cells.Item($srcRow, $col).Value2
        $ws.Cells.Item($destRow, $col).Value2 = $val
    }
}
